$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New column G width (~44.25 chars; engine rounds to nearest 1/7 px unit)
$ws.Columns.Item(7).ColumnWidth = 43.55

# G4:G5 merged note cell (left aligned, vertical centered keeps default)
$ws.Range("G4:G5").HorizontalAlignment = -4131
$ws.Range("G4").Value = "是否补充newsDesk, sectionName,差不多。"
$ws.Range("G4:G5").Merge()

# New rows 9 and 10 - addFeatureQuestion study results
$ws.Range("A9").Value = "addFeatureQuestion"
$ws.Range("B9").Value = 0.90542440000000002
$ws.Range("C9").Value = 0.92191000000000001
$ws.Range("D9").Value = 0.89119000000000004

$ws.Range("A10").Value = "addFeatureQuestion"
$ws.Range("B10").Value = 0.89829599999999998
$ws.Range("C10").Value = 0.92601
$ws.Range("D10").Value = 0.88880999999999999
$ws.Range("F10").Value = "clean newsDesk&sectionName"

$ws.Range("F9").Value = "not clean"

$ws.Range("A9:D10").Interior.Color = 5296274
$ws.Range("F9:F10").Interior.Color = 5296274

# Update selection to match
$ws.Range("F10").Select()
